# Add new column 'Event' to Card24 by admin
#
# The sheet currently spans A1:L12 with a bold/bordered header row (row 1).
# We append a new "Event" header in column M, matching the existing header
# style, and touch the rest of the column (rows 2-12) so the used range
# (dimension) grows to A1:M12, leaving those data cells blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Copy the formatting of the last existing header cell (L1, the "Date"
# header) onto the new header cell M1, then set its text.
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("M1").Value = "Event"

# Touch M2:M12 so the column becomes part of the sheet's used range while
# staying empty, matching the new blank "Event" cells for every data row.
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 13).Font.Bold = $false
}
